$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Contest 46 (row 55): MI vs DC results
$ws.Range("E55").Value = 40
$ws.Range("H55").Value = 20
$ws.Range("K55").Value = 100
$ws.Range("N55").Value = 60
$ws.Range("Q55").Value = 0
$ws.Range("T55").Value = 80

# Contest 47 (row 56): RR vs CSK results
$ws.Range("E56").Value = 80
$ws.Range("H56").Value = 100
$ws.Range("K56").Value = 60
$ws.Range("N56").Value = 20
$ws.Range("Q56").Value = 0
$ws.Range("T56").Value = 40

# Expand the season-total SUM formulas to cover the full match range (rows 10:65)
$ws.Range("E68").Formula = "=SUM(D10:D65)"
$ws.Range("H68").Formula = "=SUM(G10:G65)"
$ws.Range("K68").Formula = "=SUM(J10:J65)"
$ws.Range("N68").Formula = "=SUM(M10:M65)"
$ws.Range("Q68").Formula = "=SUM(P10:P65)"
$ws.Range("T68").Formula = "=SUM(S10:S65)"

$wb.Application.Calculate()
